$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")
$ws.Range("A13").Value = "The EU EPS uses values from the US EPS."
